# Update the absolute-path / revision metadata is managed by Excel itself on
# save and is not exposed through the object model, so we focus on the
# actual content/structure changes described by the diff: the row-2 data is
# replaced with a new course record, the mailto: hyperlink on N2 is removed
# (along with its hyperlink styling), the metacurso URL hyperlink on V2 is
# repointed to the new course, and the view/selection is reset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the N2 (director e-mail) hyperlink, keep the V2 (metacurso URL)
#     hyperlink but repoint it at the new course. Looping + deleting the
#     loop variable is what actually removes the hyperlink in this host. ---
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$N$2') {
        $hl.Delete()
    }
}
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$V$2') {
        $hl.Address = "https://utpl.instructure.com/courses/56814"
    }
}
$ws.Range("V2").Value = "https://utpl.instructure.com/courses/56814"

# N2 loses its hyperlink look (underline/theme-color "Hipervinculo" style) and
# goes back to the same plain Arial 10 formatting as the rest of row 2 - copy
# that formatting over from a neighboring plain cell.
$ws.Range("Q2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

# --- New course record values (row 2) ---
$ws.Range("A2").Value = 28
$ws.Range("E2").Value = "EDUC_7094"
$ws.Range("F2").Value = "Una Nueva Mirada a la Orientación y Asesoramiento Familiar"
$ws.Range("G2").Value = "Beltrán Guevara Patricia Maricela"
$ws.Range("H2").Value = "pmbeltran@utpl.edu.ec"
$ws.Range("I2").Value = 1103442891
$ws.Range("K2").Value = "Ciencias Sociales, Educación y Humanidades"
$ws.Range("L2").Value = "Maestría en Educación con Mención en Orientación Familiar"
$ws.Range("M2").Value = "Carrera Herrera Xiomara Paola"
$ws.Range("N2").Value = "xpcarrera@utpl.edu.ec"
$ws.Range("O2").Value = 2
$ws.Range("R2").Value = "Unidad de Formación Disciplinar Avanzada"
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = "Maestría profesional"
$ws.Range("X2").Value = "EDUC_7094_NLG_META"
$ws.Range("Y2").Value = ""
$ws.Range("AA2").Value = "EDUC_7094"
$ws.Range("AB2").Value = 45855
$ws.Range("AC2").Value = 45855

# --- Reset the view: no frozen/scrolled top-left cell, select A2:AC2 ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A2:AC2").Select()
